$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.12468433333333
$ws.Range("H2").Value = 48.374053
$ws.Range("I2").Value = 0.2955490655206278
$ws.Range("J2").Value = 0.2955490655206279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.993142333333334
$ws.Range("N2").Value = 8.979427000000001
$ws.Range("O2").Value = 0.03484385887642424
$ws.Range("P2").Value = 0.03484385887642424
$ws.Range("Q2").Value = 48.26347528973677
$ws.Range("R2").Value = 434.371277607631
$ws.Range("S2").Value = 0.01029806993005982
$ws.Range("T2").Value = 0.01029806993005982
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.12468433333333
$ws.Range("H3").Value = 48.374053
$ws.Range("I3").Value = 0.2955490655206278
$ws.Range("J3").Value = 0.2955490655206279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 31.995262
$ws.Range("N3").Value = 95.985786
$ws.Range("O3").Value = 0.3724642097459734
$ws.Range("P3").Value = 0.3724642097459735
$ws.Range("Q3").Value = 515.9134999122953
$ws.Range("R3").Value = 4643.221499210658
$ws.Range("S3").Value = 0.1100814491303016
$ws.Range("T3").Value = 0.1100814491303016
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.12468433333333
$ws.Range("H4").Value = 48.374053
$ws.Range("I4").Value = 0.2955490655206278
$ws.Range("J4").Value = 0.2955490655206279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 34.28929533333334
$ws.Range("N4").Value = 102.867886
$ws.Range("O4").Value = 0.3991695798295478
$ws.Range("P4").Value = 0.3991695798295478
$ws.Range("Q4").Value = 552.9040632624398
$ws.Range("R4").Value = 4976.136569361957
$ws.Range("S4").Value = 0.1179741963028845
$ws.Range("T4").Value = 0.1179741963028845
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.12468433333333
$ws.Range("H5").Value = 48.374053
$ws.Range("I5").Value = 0.2955490655206278
$ws.Range("J5").Value = 0.2955490655206279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.62387466666667
$ws.Range("N5").Value = 49.871624
$ws.Range("O5").Value = 0.1935223515480544
$ws.Range("P5").Value = 0.1935223515480545
$ws.Range("Q5").Value = 268.0547313968968
$ws.Range("R5").Value = 2412.492582572072
$ws.Range("S5").Value = 0.05719535015738191
$ws.Range("T5").Value = 0.05719535015738193
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.68975
$ws.Range("H6").Value = 44.06925
$ws.Range("I6").Value = 0.2692481784748309
$ws.Range("J6").Value = 0.2692481784748309
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.993142333333334
$ws.Range("N6").Value = 8.979427000000001
$ws.Range("O6").Value = 0.03484385887642424
$ws.Range("P6").Value = 0.03484385887642424
$ws.Range("Q6").Value = 43.96851259108333
$ws.Range("R6").Value = 395.71661331975
$ws.Range("S6").Value = 0.009381645533511293
$ws.Range("T6").Value = 0.009381645533511295
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.68975
$ws.Range("H7").Value = 44.06925
$ws.Range("I7").Value = 0.2692481784748309
$ws.Range("J7").Value = 0.2692481784748309
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.995262
$ws.Range("N7").Value = 95.985786
$ws.Range("O7").Value = 0.3724642097459734
$ws.Range("P7").Value = 0.3724642097459735
$ws.Range("Q7").Value = 470.0023999644999
$ws.Range("R7").Value = 4230.0215996805
$ws.Range("S7").Value = 0.1002853100211707
$ws.Range("T7").Value = 0.1002853100211707
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.68975
$ws.Range("H8").Value = 44.06925
$ws.Range("I8").Value = 0.2692481784748309
$ws.Range("J8").Value = 0.2692481784748309
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.28929533333334
$ws.Range("N8").Value = 102.867886
$ws.Range("O8").Value = 0.3991695798295478
$ws.Range("P8").Value = 0.3991695798295478
$ws.Range("Q8").Value = 503.7011761228333
$ws.Range("R8").Value = 4533.3105851055
$ws.Range("S8").Value = 0.1074756822716693
$ws.Range("T8").Value = 0.1074756822716693
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.68975
$ws.Range("H9").Value = 44.06925
$ws.Range("I9").Value = 0.2692481784748309
$ws.Range("J9").Value = 0.2692481784748309
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 16.62387466666667
$ws.Range("N9").Value = 49.871624
$ws.Range("O9").Value = 0.1935223515480544
$ws.Range("P9").Value = 0.1935223515480545
$ws.Range("Q9").Value = 244.2005628846666
$ws.Range("R9").Value = 2197.805065962
$ws.Range("S9").Value = 0.05210554064847953
$ws.Range("T9").Value = 0.05210554064847954
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.967860666666667
$ws.Range("H10").Value = 5.903582
$ws.Range("I10").Value = 0.03606888476606249
$ws.Range("J10").Value = 0.03606888476606249
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.993142333333334
$ws.Range("N10").Value = 8.979427000000001
$ws.Range("O10").Value = 0.03484385887642424
$ws.Range("P10").Value = 0.03484385887642424
$ws.Range("Q10").Value = 5.890087067501556
$ws.Range("R10").Value = 53.01078360751401
$ws.Range("S10").Value = 0.001256779130618689
$ws.Range("T10").Value = 0.00125677913061869
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.967860666666667
$ws.Range("H11").Value = 5.903582
$ws.Range("I11").Value = 0.03606888476606249
$ws.Range("J11").Value = 0.03606888476606249
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 31.995262
$ws.Range("N11").Value = 95.985786
$ws.Range("O11").Value = 0.3724642097459734
$ws.Range("P11").Value = 0.3724642097459735
$ws.Range("Q11").Value = 62.96221760949467
$ws.Range("R11").Value = 566.6599584854521
$ws.Range("S11").Value = 0.01343436866081004
$ws.Range("T11").Value = 0.01343436866081005
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.967860666666667
$ws.Range("H12").Value = 5.903582
$ws.Range("I12").Value = 0.03606888476606249
$ws.Range("J12").Value = 0.03606888476606249
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 34.28929533333334
$ws.Range("N12").Value = 102.867886
$ws.Range("O12").Value = 0.3991695798295478
$ws.Range("P12").Value = 0.3991695798295478
$ws.Range("Q12").Value = 67.47655557418356
$ws.Range("R12").Value = 607.289000167652
$ws.Range("S12").Value = 0.01439760157698954
$ws.Range("T12").Value = 0.01439760157698954
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.967860666666667
$ws.Range("H13").Value = 5.903582
$ws.Range("I13").Value = 0.03606888476606249
$ws.Range("J13").Value = 0.03606888476606249
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.62387466666667
$ws.Range("N13").Value = 49.871624
$ws.Range("O13").Value = 0.1935223515480544
$ws.Range("P13").Value = 0.1935223515480545
$ws.Range("Q13").Value = 32.71346908412978
$ws.Range("R13").Value = 294.421221757168
$ws.Range("S13").Value = 0.00698013539764421
$ws.Range("T13").Value = 0.006980135397644211
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 21.77610566666667
$ws.Range("H14").Value = 65.328317
$ws.Range("I14").Value = 0.3991338712384788
$ws.Range("J14").Value = 0.3991338712384788
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.993142333333334
$ws.Range("N14").Value = 8.979427000000001
$ws.Range("O14").Value = 0.03484385887642424
$ws.Range("P14").Value = 0.03484385887642424
$ws.Range("Q14").Value = 65.17898372603989
$ws.Range("R14").Value = 586.6108535343591
$ws.Range("S14").Value = 0.01390736428223444
$ws.Range("T14").Value = 0.01390736428223444
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 21.77610566666667
$ws.Range("H15").Value = 65.328317
$ws.Range("I15").Value = 0.3991338712384788
$ws.Range("J15").Value = 0.3991338712384788
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 31.995262
$ws.Range("N15").Value = 95.985786
$ws.Range("O15").Value = 0.3724642097459734
$ws.Range("P15").Value = 0.3724642097459735
$ws.Range("Q15").Value = 696.7322061446847
$ws.Range("R15").Value = 6270.589855302163
$ws.Range("S15").Value = 0.1486630819336911
$ws.Range("T15").Value = 0.1486630819336911
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 21.77610566666667
$ws.Range("H16").Value = 65.328317
$ws.Range("I16").Value = 0.3991338712384788
$ws.Range("J16").Value = 0.3991338712384788
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 34.28929533333334
$ws.Range("N16").Value = 102.867886
$ws.Range("O16").Value = 0.3991695798295478
$ws.Range("P16").Value = 0.3991695798295478
$ws.Range("Q16").Value = 746.6873184142069
$ws.Range("R16").Value = 6720.185865727862
$ws.Range("S16").Value = 0.1593220996780044
$ws.Range("T16").Value = 0.1593220996780044
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 21.77610566666667
$ws.Range("H17").Value = 65.328317
$ws.Range("I17").Value = 0.3991338712384788
$ws.Range("J17").Value = 0.3991338712384788
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 16.62387466666667
$ws.Range("N17").Value = 49.871624
$ws.Range("O17").Value = 0.1935223515480544
$ws.Range("P17").Value = 0.1935223515480545
$ws.Range("Q17").Value = 362.0032513307564
$ws.Range("R17").Value = 3258.029261976808
$ws.Range("S17").Value = 0.07724132534454878
$ws.Range("T17").Value = 0.07724132534454879
